$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 597.98
$ws.Range("J17").Value = 599.9798
$ws.Range("L17").Value = 1799.9394
$ws.Range("N17").Value = -2135.9394
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents() | Out-Null
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents() | Out-Null
$ws.Range("H92").Value = 356.8421
$ws.Range("I92").Value = 280.76923
$ws.Range("J92").Value = 521.6667
$ws.Range("K92").Value = 280.76923
$ws.Range("L92").Value = 521.6667
$ws.Range("M92").Value = 967.23077
$ws.Range("N92").Value = -3017.6667
$ws.Range("H137").Value = 1484.4062
$ws.Range("I137").Value = 933.7778
$ws.Range("J137").Value = 2192.3572
$ws.Range("K137").Value = 2801.3334
$ws.Range("L137").Value = 6577.071599999999
$ws.Range("M137").Value = -251.3334
$ws.Range("N137").Value = -11677.0716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5596.7
$ws.Range("I32").Value = 4496.3335
$ws.Range("J32").Value = 15500
$ws.Range("K32").Value = 4496.3335
$ws.Range("L32").Value = 15500
$ws.Range("M32").Value = -4209.3335
$ws.Range("N32").Value = -16074
$ws.Range("H61").Value = 2140.543
$ws.Range("I61").Value = 1836.3928
$ws.Range("K61").Value = 1836.3928
$ws.Range("M61").Value = -1624.3928
$ws.Range("H74").Value = 12821878
$ws.Range("I74").Value = 15152664
$ws.Range("J74").Value = 2557.5
$ws.Range("K74").Value = 15152664
$ws.Range("L74").Value = 2557.5
$ws.Range("M74").Value = -15151790
$ws.Range("N74").Value = -4305.5
$ws.Range("H77").Value = 12821878
$ws.Range("I77").Value = 15152664
$ws.Range("J77").Value = 2557.5
$ws.Range("K77").Value = 75763320
$ws.Range("L77").Value = 12787.5
$ws.Range("M77").Value = -75758952
$ws.Range("N77").Value = -21523.5
$ws.Range("H88").Value = 2500
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -2094
$ws.Range("N88").ClearContents() | Out-Null
$ws.Range("H91").Value = 2500
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -1096
$ws.Range("N91").ClearContents() | Out-Null
$ws.Range("H136").Value = 2140.543
$ws.Range("I136").Value = 1836.3928
$ws.Range("K136").Value = 5509.178400000001
$ws.Range("M136").Value = -2959.178400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1308.4193
$ws.Range("I86").Value = 1279.7778
$ws.Range("J86").Value = 1501.75
$ws.Range("K86").Value = 1279.7778
$ws.Range("L86").Value = 1501.75
$ws.Range("M86").Value = -156.7778000000001
$ws.Range("N86").Value = -3747.75
$ws.Range("H89").Value = 1308.4193
$ws.Range("I89").Value = 1279.7778
$ws.Range("J89").Value = 1501.75
$ws.Range("K89").Value = 6398.889
$ws.Range("L89").Value = 7508.75
$ws.Range("M89").Value = -782.8890000000001
$ws.Range("N89").Value = -18740.75
$ws.Range("H105").Value = 1340310
$ws.Range("I105").Value = 2068033.8
$ws.Range("J105").Value = 6150
$ws.Range("K105").Value = 2068033.8
$ws.Range("L105").Value = 6150
$ws.Range("M105").Value = -2066286.8
$ws.Range("N105").Value = -9644
$ws.Range("H134").Value = 4099778.2
$ws.Range("I134").Value = 5556591
$ws.Range("J134").Value = 2491.5625
$ws.Range("K134").Value = 16669773
$ws.Range("L134").Value = 7474.6875
$ws.Range("M134").Value = -16667238
$ws.Range("N134").Value = -12544.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10641713
$ws.Range("I31").Value = 21740258
$ws.Range("J31").Value = 5605.9165
$ws.Range("K31").Value = 21740258
$ws.Range("L31").Value = 5605.9165
$ws.Range("M31").Value = -21739963
$ws.Range("N31").Value = -6195.9165
$ws.Range("H34").Value = 10641713
$ws.Range("I34").Value = 21740258
$ws.Range("J34").Value = 5605.9165
$ws.Range("K34").Value = 21740258
$ws.Range("L34").Value = 5605.9165
$ws.Range("M34").Value = -21740056
$ws.Range("N34").Value = -6009.9165
$ws.Range("H122").Value = 6287.625
$ws.Range("I122").Value = 7347.2354
$ws.Range("J122").Value = 3714.2856
$ws.Range("K122").Value = 22041.7062
$ws.Range("L122").Value = 11142.8568
$ws.Range("M122").Value = -19591.7062
$ws.Range("N122").Value = -16042.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4345.136
$ws.Range("I3").Value = 3508.5715
$ws.Range("J3").Value = 4735.533
$ws.Range("K3").Value = 10525.7145
$ws.Range("L3").Value = 14206.599
$ws.Range("M3").Value = -10413.7145
$ws.Range("N3").Value = -14430.599
$ws.Range("H60").Value = 573.6842
$ws.Range("I60").Value = 265.33334
$ws.Range("J60").Value = 1730
$ws.Range("K60").Value = 796.0000200000001
$ws.Range("L60").Value = 5190
$ws.Range("M60").Value = -545.0000200000001
$ws.Range("N60").Value = -5692
$ws.Range("H131").Value = 919.71875
$ws.Range("I131").Value = 316.66666
$ws.Range("J131").Value = 982.10345
$ws.Range("K131").Value = 949.9999799999999
$ws.Range("L131").Value = 2946.31035
$ws.Range("M131").Value = 4090.00002
$ws.Range("N131").Value = -13026.31035

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 55000
$ws.Range("J68").Value = 55000
$ws.Range("L68").Value = 55000
$ws.Range("N68").Value = -56622
$ws.Range("H71").Value = 55000
$ws.Range("J71").Value = 55000
$ws.Range("L71").Value = 165000
$ws.Range("N71").Value = -173112

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 523.4
$ws.Range("I22").Value = 223
$ws.Range("J22").Value = 623.5333
$ws.Range("K22").Value = 223
$ws.Range("L22").Value = 623.5333
$ws.Range("M22").Value = 72
$ws.Range("N22").Value = -1213.5333
$ws.Range("H27").Value = 523.4
$ws.Range("I27").Value = 223
$ws.Range("J27").Value = 623.5333
$ws.Range("K27").Value = 223
$ws.Range("L27").Value = 623.5333
$ws.Range("M27").Value = -116
$ws.Range("N27").Value = -837.5333
$ws.Range("H40").Value = 2415.3447
$ws.Range("I40").Value = 1625
$ws.Range("J40").Value = 2973.2354
$ws.Range("K40").Value = 1625
$ws.Range("L40").Value = 2973.2354
$ws.Range("M40").Value = -1489
$ws.Range("N40").Value = -3245.2354
$ws.Range("H96").Value = 29500
$ws.Range("J96").Value = 29500
$ws.Range("L96").Value = 29500
$ws.Range("N96").Value = -34992
$ws.Range("H122").Value = 3111.2812
$ws.Range("I122").Value = 2944.6538
$ws.Range("J122").Value = 3833.3333
$ws.Range("K122").Value = 8833.9614
$ws.Range("L122").Value = 11499.9999
$ws.Range("M122").Value = -6383.9614
$ws.Range("N122").Value = -16399.9999
$ws.Range("H136").Value = 2612.7812
$ws.Range("I136").Value = 1272.0952
$ws.Range("K136").Value = 3816.2856
$ws.Range("M136").Value = -1266.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1322.2424
$ws.Range("I122").Value = 834.2727
$ws.Range("J122").Value = 2298.182
$ws.Range("K122").Value = 2502.8181
$ws.Range("L122").Value = 6894.545999999999
$ws.Range("M122").Value = -52.81809999999996
$ws.Range("N122").Value = -11794.546
$ws.Range("H136").Value = 3262.4546
$ws.Range("I136").Value = 996.46155
$ws.Range("J136").Value = 8785.8125
$ws.Range("K136").Value = 2989.38465
$ws.Range("L136").Value = 26357.4375
$ws.Range("M136").Value = -439.38465
$ws.Range("N136").Value = -31457.4375
